# Update the financial figures on the "company_list" sheet so that the
# (previously mis-scaled) raw values are replaced with the corrected
# (smaller-magnitude) figures, and drop the stray extra forecast-year
# rows (7, 8, 9) down to just their id / label columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected values for rows 2-6 (columns D..AJ) -------------------------
$newValues = @{
    "D2" = 6032
    "E2" = 65
    "F2" = 65
    "G2" = 45
    "H2" = 50
    "I2" = 49
    "J2" = 1
    "K2" = 7070
    "L2" = 4766
    "M2" = 2304
    "N2" = 2251
    "O2" = 53
    "P2" = 237
    "Q2" = -295
    "R2" = -252
    "S2" = 550
    "T2" = 284
    "U2" = -579
    "V2" = 2264
    "W2" = 1.07
    "X2" = 0.83
    "Y2" = 2.21
    "Z2" = 0.78
    "AA2" = 206.87
    "AB2" = 871.52
    "AC2" = 207
    "AD2" = 47.54
    "AE2" = 10861
    "AF2" = 0.91
    "AG2" = 60
    "AH2" = 0.61
    "AI2" = 25.3
    "AJ2" = 23728210

    "D3" = 5835
    "E3" = 54
    "F3" = 54
    "G3" = 50
    "H3" = 35
    "I3" = 33
    "J3" = 2
    "K3" = 6883
    "L3" = 4552
    "M3" = 2331
    "N3" = 2278
    "O3" = 53
    "P3" = 237
    "Q3" = 154
    "R3" = -190
    "S3" = 35
    "T3" = 177
    "U3" = -23
    "V3" = 2405
    "W3" = 0.93
    "X3" = 0.61
    "Y3" = 1.47
    "Z3" = 0.51
    "AA3" = 195.25
    "AB3" = 876.3200000000001
    "AC3" = 141
    "AD3" = 65.25
    "AE3" = 10989
    "AF3" = 0.83
    "AG3" = 50
    "AH3" = 0.55
    "AI3" = 31.08
    "AJ3" = 23728210

    "D4" = 5854
    "E4" = -130
    "F4" = -130
    "G4" = -182
    "H4" = -180
    "I4" = -179
    "J4" = -1
    "K4" = 6672
    "L4" = 4521
    "M4" = 2151
    "N4" = 2100
    "O4" = 52
    "P4" = 237
    "Q4" = 3
    "R4" = -114
    "S4" = 70
    "T4" = 133
    "U4" = -129
    "V4" = 2605
    "W4" = -2.22
    "X4" = -3.07
    "Y4" = -8.19
    "Z4" = -2.65
    "AA4" = 210.12
    "AB4" = 801.21
    "AC4" = -755
    "AD4" = -9.470000000000001
    "AE4" = 10130
    "AF4" = 0.71
    "AG4" = 40
    "AH4" = 0.5600000000000001
    "AI4" = -4.63
    "AJ4" = 23728210

    "D5" = 6101
    "E5" = 172
    "F5" = 172
    "G5" = 135
    "H5" = 99
    "I5" = 96
    "J5" = 2
    "K5" = 7306
    "L5" = 5117
    "M5" = 2189
    "N5" = 2133
    "O5" = 56
    "P5" = 237
    "Q5" = 106
    "R5" = -191
    "S5" = 102
    "T5" = 147
    "U5" = -40
    "V5" = 2843
    "W5" = 2.82
    "X5" = 1.62
    "Y5" = 4.54
    "Z5" = 1.41
    "AA5" = 233.78
    "AB5" = 828.21
    "AC5" = 405
    "AD5" = 23.68
    "AE5" = 10292
    "AF5" = 0.93
    "AG5" = 60
    "AH5" = 0.63
    "AI5" = 12.93
    "AJ5" = 23728210

    "D6" = 6548
    "E6" = 52
    "F6" = 52
    "G6" = 18
    "H6" = 15
    "I6" = 12
    "K6" = 7727
    "L6" = 5552
    "M6" = 2175
    "N6" = 2119
    "P6" = 237
    "Q6" = 39
    "R6" = -503
    "S6" = 524
    "T6" = 568
    "U6" = -528
    "V6" = 3474
    "W6" = 0.8
    "X6" = 0.22
    "Y6" = 0.58
    "Z6" = 0.19
    "AA6" = 255.27
    "AB6" = 820.62
    "AC6" = 52
    "AD6" = 99.90000000000001
    "AE6" = 10224
    "AF6" = 0.51
    "AG6" = 50
    "AH6" = 0.97
    "AI6" = 84.23999999999999
    "AJ6" = 23728210
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

# --- Rows 7-9: keep only the id (A) / period (B) / label (C) columns -------
# these rows previously carried a full set of (erroneous) figures; the fix
# drops all of the numeric/metric columns, leaving the row identifying info.
$dataCols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")

foreach ($r in 7..9) {
    foreach ($c in $dataCols) {
        $ws.Range("$c$r").ClearContents()
    }
}
